$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at row 2 for "Pollo Grangero" (shifts existing rows down by one)
$ws.Rows.Item(2).Insert()

# Copy number formatting (date style) from row 3 (agroplus, shifted down) into new row 2
# so that G2/H2 keep the same date style used by the rest of the createdAt/updatedAt columns.
$ws.Range("G3:H3").Copy()
$ws.Range("G2:H2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new "Pollo Grangero" row
$ws.Range("A2").Value = "67bf38f0a9c0f101fb8c8295"
$ws.Range("B2").Value = "Pollo Grangero"
$ws.Range("C2").Value = "Medio"
$ws.Range("D2").Value = 2010
$ws.Range("E2").Value = "Food"
$ws.Range("F2").Value = "Empresa de comida."
$ws.Range("G2").Value = 45714.41199854167
$ws.Range("H2").Value = 45714.430585925926

# 2) Remove the "campero" row, which is now at row 5 after the insertion above
$ws.Rows.Item(5).Delete()

# 3) Populate the new "_id" column (A) for the remaining companies that shifted down by one row
$ws.Range("A3").Value = "67bf3963a9c0f101fb8c82b1"  # agroplus
$ws.Range("A4").Value = "67bf3931a9c0f101fb8c82a5"  # autoworld
$ws.Range("A5").Value = "67bf3911a9c0f101fb8c829d"  # ecogreen
$ws.Range("A6").Value = "67bf3957a9c0f101fb8c82ad"  # edusmart
$ws.Range("A7").Value = "67bf394ca9c0f101fb8c82a9"  # fintrust
$ws.Range("A8").Value = "67bf3926a9c0f101fb8c82a1"  # medicare
$ws.Range("A9").Value = "67bf3901a9c0f101fb8c8299"  # techsoft
